$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E (4:9) takes on the values that used to live in column I,
# then column I is cleared/removed entirely.
$ws.Range("E4").Value = 195
$ws.Range("E5").Value = 40
$ws.Range("E6").Value = 60
$ws.Range("E7").Value = 30
$ws.Range("E8").Value = 50
$ws.Range("E9").Value = 20

# Remove the now-redundant column I data (I4:I9).
$ws.Range("I4:I9").ClearContents()

# E11 keeps its SUM formula; it will recalc to 395 once E4:E9 change.
$ws.Range("E11").Formula = "=SUM(E4:E9)"

# Update the "marc" price and the final total.
$ws.Range("E13").Value = 250
$ws.Range("E17").Value = 795

# Move the active selection to H7, matching the saved view state.
$ws.Range("H7").Select()
